# edit.ps1
# Implements the commit: add a new "aws.ses" command category to the
# '#system' reference sheet (with sendMail/sendTextMail commands), add a
# new "base64(var,file)" command to the "io" category, and add a new
# "upload(url,body,fileParams,var)" command to the "ws" category. Also
# updates every defined name so that it still points at the right range
# after the new column/rows are inserted.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1. Insert a brand-new column before the old column C ("base"). This
#    shifts every column from C..Z one position to the right (D..AA) but
#    leaves columns A (target) and B (aws.s3) untouched.
# ---------------------------------------------------------------------
$ws.Columns("C:C").Insert()

# Populate the new column C with the "aws.ses" category header plus its
# two commands.
$ws.Cells.Item(1, 3).Value2 = "aws.ses"
$ws.Cells.Item(2, 3).Value2 = "sendMail(profile,to,subject,body)"
$ws.Cells.Item(3, 3).Value2 = "sendTextMail(profile,to,subject,body)"

# ---------------------------------------------------------------------
# 2. Insert "aws.ses" into the "target" list (column A), right after
#    "aws.s3" (row 2), shifting the remaining category names down one
#    row (old A3:A26 -> A4:A27). Doing this with manual cell copies
#    (instead of Range.Insert) keeps the shift confined to column A.
# ---------------------------------------------------------------------
for ($r = 26; $r -ge 3; $r--) {
    $ws.Cells.Item($r + 1, 1).Value2 = $ws.Cells.Item($r, 1).Value2
}
$ws.Cells.Item(3, 1).Value2 = "aws.ses"

# ---------------------------------------------------------------------
# 3. Insert "base64(var,file)" into the "io" list, which now lives in
#    column J (was column I before the column insert above). It is
#    inserted alphabetically at row 5 (between "assertReadableFile..."
#    and "compare..."), shifting the rest of column J down one row
#    (old J5:J23 -> J6:J24).
# ---------------------------------------------------------------------
for ($r = 23; $r -ge 5; $r--) {
    $ws.Cells.Item($r + 1, 10).Value2 = $ws.Cells.Item($r, 10).Value2
}
$ws.Cells.Item(5, 10).Value2 = "base64(var,file)"

# ---------------------------------------------------------------------
# 4. Append "upload(url,body,fileParams,var)" to the "ws" list, which
#    now lives in column Y (was column X before the column insert
#    above). Alphabetically it sorts after "soap(...)", i.e. at the end
#    of the existing list, so it is simply appended as a new last row.
# ---------------------------------------------------------------------
$ws.Cells.Item(17, 25).Value2 = "upload(url,body,fileParams,var)"

# ---------------------------------------------------------------------
# 5. Fix up all the defined names so they still point at the correct
#    (now shifted) ranges, and add the new "aws.ses" named range.
# ---------------------------------------------------------------------
$sys = "'#system'"

$wb.Names.Item("base").RefersTo        = "=$sys!`$D`$2:`$D`$36"
$wb.Names.Item("csv").RefersTo         = "=$sys!`$E`$2:`$E`$5"
$wb.Names.Item("desktop").RefersTo     = "=$sys!`$F`$2:`$F`$92"
$wb.Names.Item("excel").RefersTo       = "=$sys!`$G`$2:`$G`$14"
$wb.Names.Item("external").RefersTo    = "=$sys!`$H`$2:`$H`$3"
$wb.Names.Item("image").RefersTo       = "=$sys!`$I`$2:`$I`$5"
$wb.Names.Item("io").RefersTo          = "=$sys!`$J`$2:`$J`$24"
$wb.Names.Item("jms").RefersTo         = "=$sys!`$K`$2:`$K`$4"
$wb.Names.Item("json").RefersTo        = "=$sys!`$L`$2:`$L`$14"
$wb.Names.Item("mail").RefersTo        = "=$sys!`$M`$2:`$M`$2"
$wb.Names.Item("number").RefersTo      = "=$sys!`$N`$2:`$N`$15"
$wb.Names.Item("pdf").RefersTo         = "=$sys!`$O`$2:`$O`$16"
$wb.Names.Item("rdbms").RefersTo       = "=$sys!`$P`$2:`$P`$7"
$wb.Names.Item("redis").RefersTo       = "=$sys!`$Q`$2:`$Q`$10"
$wb.Names.Item("ssh").RefersTo         = "=$sys!`$T`$2:`$T`$9"
$wb.Names.Item("step").RefersTo        = "=$sys!`$U`$2:`$U`$4"
$wb.Names.Item("target").RefersTo      = "=$sys!`$A`$2:`$A`$27"
$wb.Names.Item("web").RefersTo         = "=$sys!`$V`$2:`$V`$117"
$wb.Names.Item("webalert").RefersTo    = "=$sys!`$W`$2:`$W`$8"
$wb.Names.Item("webcookie").RefersTo   = "=$sys!`$X`$2:`$X`$8"
$wb.Names.Item("ws").RefersTo          = "=$sys!`$Y`$2:`$Y`$17"
$wb.Names.Item("xml").RefersTo         = "=$sys!`$AA`$2:`$AA`$11"
$wb.Names.Item("sms").RefersTo         = "=$sys!`$R`$2:`$R`$2"
$wb.Names.Item("sound").RefersTo       = "=$sys!`$S`$2:`$S`$5"
$wb.Names.Item("ws.async").RefersTo    = "=$sys!`$Z`$2:`$Z`$8"

# aws.s3, date, db, math, mq and nextgen were not affected by the column
# insert, so they are left as-is.

# Brand new named range for the new "aws.ses" column.
$wb.Names.Add("aws.ses", "=$sys!`$C`$2:`$C`$3")
